$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new columns before column B (shifts existing B:V data right to K:AE)
$ws.Columns("B:J").Insert()

# New date headers for row 1 (most-recent-first ordering), newest in B1 .. oldest of the new batch in J1
$newHeaders = @("Sep_08","Aug_25","Aug_04","Jul_23","Jul_17","Jul_07","Jun_30","Jun_24","Jun_16")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $newHeaders[$i]
}

# Fill the newly inserted data columns (B:J) for every data row (2-33) with the "UN" placeholder,
# matching the existing rating-unchanged marker used throughout the sheet.
for ($r = 2; $r -le 33; $r++) {
    for ($c = 2; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = "UN"
    }
}
